$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T for the 12 data rows (rows 2-13). Sending cluster now
# includes "ECs" in addition to the original "FAPs"/"sCs", per Dr Hou's
# advice, so every sending/target cluster combination is recomputed.
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rows = @(
    @("ECs", "Vtn", "Itga5", "ECs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 206.0215709003713, 1854.194138103342, 0.03479986429231236, 0.03479986429231236),
    @("ECs", "Vtn", "Itga5", "FAPs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 266.8261266495293, 2401.435139845764, 0.04507058632970669, 0.04507058632970668),
    @("ECs", "Vtn", "Itga5", "M2", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 167.3668472990147, 1506.301625691132, 0.02827055219307319, 0.02827055219307319),
    @("ECs", "Vtn", "Itga5", "sCs", 2, 0.6666666666666666, 6.597131, 19.791393, 0.1209543635982448, 0.1209543635982448, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 75.85744285909367, 682.716985731843, 0.01281336078315254, 0.01281336078315254),
    @("FAPs", "Vtn", "Itga5", "ECs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 657.4688061478755, 5917.21925533088, 0.1110554838038731, 0.1110554838038731),
    @("FAPs", "Vtn", "Itga5", "FAPs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 851.5120730836622, 7663.60865775296, 0.1438320485426587, 0.1438320485426587),
    @("FAPs", "Vtn", "Itga5", "M2", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 534.1114563951645, 4807.00310755648, 0.09021873835002164, 0.09021873835002164),
    @("FAPs", "Vtn", "Itga5", "sCs", 3, 1, 21.05317333333333, 63.15952, 0.3859970617919927, 0.3859970617919927, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 242.0809732497245, 2178.72875924752, 0.04089079109543925, 0.04089079109543925),
    @("sCs", "Vtn", "Itga5", "ECs", 3, 1, 26.89201066666667, 80.67603200000001, 0.4930485746097625, 0.4930485746097625, 3, 1, 31.22896466666667, 93.686894, 0.2877106972998646, 0.2877106972998646, 839.8096509249565, 7558.286858324608, 0.1418553492036791, 0.1418553492036791),
    @("sCs", "Vtn", "Itga5", "FAPs", 3, 1, 26.89201066666667, 80.67603200000001, 0.4930485746097625, 0.4930485746097625, 3, 1, 40.44578266666667, 121.337348, 0.3726247238124506, 0.3726247238124505, 1087.668418893682, 9789.015770043135, 0.1837220889400852, 0.1837220889400852),
    @("sCs", "Vtn", "Itga5", "M2", 3, 1, 26.89201066666667, 80.67603200000001, 0.4930485746097625, 0.4930485746097625, 3, 1, 25.36964133333333, 76.108924, 0.2337290805561598, 0.2337290805561598, 682.2406653455076, 6140.165988109568, 0.1152397900130649, 0.1152397900130649),
    @("sCs", "Vtn", "Itga5", "sCs", 3, 1, 26.89201066666667, 80.67603200000001, 0.4930485746097625, 0.4930485746097625, 3, 1, 11.49855033333333, 34.495651, 0.1059354983315251, 0.1059354983315251, 309.2191382152036, 2782.972243936832, 0.05223134645293334, 0.05223134645293334)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowData = $rows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $columns.Count; $c++) {
        $ws.Range("$($columns[$c])$r").Value = $rowData[$c]
    }
}
